# feat: add 2022-Q4 data
#
# The workbook currently has two sheets:
#   1. 总计        - summary table (one row per quarter)
#   2. 2021-Q2     - per-fund detail table for the 2021-Q2 quarter
#
# We add a new quarter "2022-Q4":
#   - a new row is added to 总计 for 2022-Q4 (placed right after the header,
#     before the existing 2021-Q2 summary row, which shifts down one row)
#   - a brand-new sheet "2022-Q4" is inserted right after 总计 (and before
#     2021-Q2) holding the per-fund detail rows for that quarter

$wb = $excel.ActiveWorkbook
$totalWs = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计" (i.e. before "2021-Q2")
# ---------------------------------------------------------------------------
$q4Ws = $wb.Worksheets.Add($null, $totalWs)
$q4Ws.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert the 2022-Q4 summary row just under the
#    header, pushing the existing 2021-Q2 summary row down to row 3.
# ---------------------------------------------------------------------------

# Row 3 <- old row 2 content ("2021-Q2", 4, 0.32), index column value 1.
$totalWs.Range("A2").Copy($totalWs.Range("A3"))
$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2021-Q2"
$totalWs.Range("C3").Value = 4
$totalWs.Range("D3").Value = 0.32

# Row 2 <- new 2022-Q4 summary data, index column value 0 (already there).
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 2
$totalWs.Range("D2").Value = 0.02

# ---------------------------------------------------------------------------
# 3. Populate the new "2022-Q4" detail sheet.
# ---------------------------------------------------------------------------

# Header row (B1:H1), formatted the same way as the "总计" header cells.
$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q4Cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $q4Cols.Length; $i++) {
    $cell = $q4Ws.Range($q4Cols[$i] + "1")
    $totalWs.Range("B1").Copy($cell)
    $cell.Value = $q4Headers[$i]
}

# Index column (A2:A3), formatted the same way as "总计"'s index column.
$totalWs.Range("A2").Copy($q4Ws.Range("A2"))
$q4Ws.Range("A2").Value = 0
$totalWs.Range("A2").Copy($q4Ws.Range("A3"))
$q4Ws.Range("A3").Value = 1

# Data row 2: 002236 / 大成中证360互联网+大数据100指数A
$q4Ws.Range("B2").NumberFormat = "@"
$q4Ws.Range("B2").Value = "002236"
$q4Ws.Range("C2").Value = "大成中证360互联网+大数据100指数A"
$q4Ws.Range("D2").NumberFormat = "@"
$q4Ws.Range("D2").Value = "1.15"
$q4Ws.Range("E2").NumberFormat = "@"
$q4Ws.Range("E2").Value = "92.50"
$q4Ws.Range("F2").NumberFormat = "@"
$q4Ws.Range("F2").Value = "1.00"
$q4Ws.Range("G2").NumberFormat = "@"
$q4Ws.Range("G2").Value = "0.0115"
$q4Ws.Range("H2").Value = 7

# Data row 3: 003359 / 大成中证360互联网+大数据100指数C
$q4Ws.Range("B3").NumberFormat = "@"
$q4Ws.Range("B3").Value = "003359"
$q4Ws.Range("C3").Value = "大成中证360互联网+大数据100指数C"
$q4Ws.Range("D3").NumberFormat = "@"
$q4Ws.Range("D3").Value = "1.12"
$q4Ws.Range("E3").NumberFormat = "@"
$q4Ws.Range("E3").Value = "92.50"
$q4Ws.Range("F3").NumberFormat = "@"
$q4Ws.Range("F3").Value = "1.00"
$q4Ws.Range("G3").NumberFormat = "@"
$q4Ws.Range("G3").Value = "0.0112"
$q4Ws.Range("H3").Value = 7

Write-Output "2022-Q4 sheet and summary row added"
